$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 837, shifting existing rows 837:878 down to 838:879
$ws.Rows.Item(837).Insert()

# Seed the new row 837 by duplicating the row above (same date/weekday/ranking,
# all stored as plain text like the rest of the sheet) then fix the time value.
$ws.Range("A836:D836").Copy()
$ws.Range("A837").PasteSpecial()
$ws.Range("C837").Value = 16
